$wb = $excel.ActiveWorkbook

# Add a new worksheet after the last existing sheet and name it "Sheet2"
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws2.Name = "Sheet2"

# Populate the new sheet with values
$ws2.Range("A1").Value = "multiple"
$ws2.Range("A2").Value = "sheets"
$ws2.Range("A3").Value = "test"

# Select A3 and make Sheet2 the active sheet/tab
$ws2.Range("A3").Select()
$ws2.Activate()
